$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "RM 232" (row 26) entirely
$ws.Rows.Item(26).Delete()

# Delete the row for "SC 92" (was row 28, now row 27 after the previous delete) entirely
$ws.Rows.Item(27).Delete()

# Apply remaining value changes (row numbers are now in the final, post-deletion layout)
$ws.Range("D6").Value = -14.2      # RM 21
$ws.Range("D8").Value = ""         # RM 38
$ws.Range("D19").Value = -15.5     # RM 125
$ws.Range("D21").Value = ""        # RM 135
$ws.Range("D23").Value = -13.9     # RM 140
$ws.Range("C26").Value = ""        # SC 5
$ws.Range("C27").Value = 10        # SC 101
$ws.Range("D27").Value = ""        # SC 101
$ws.Range("C29").Value = ""        # SC 119
$ws.Range("D29").Value = -13       # SC 119

Write-Host "Done"
